$wb = $excel.ActiveWorkbook

# --- Append newly logged play-by-play values to running string logs (Week 13) ---
$ws_YDS = $wb.Worksheets.Item("YDS")
$ws_ST = $wb.Worksheets.Item("ST")
$ws_YDS.Range("B2").Value = $ws_YDS.Range("B2").Value2 + " 2 1 3 17 9 1 0 1 15 2 2 3 2 3 8 17"
$ws_YDS.Range("B3").Value = $ws_YDS.Range("B3").Value2 + " 8 7 3 10 7 5 10 29 4 19 5 18 14 19 11 -2 15 4 14 7 11"
$ws_YDS.Range("C2").Value = $ws_YDS.Range("C2").Value2 + " 0 5 9 2 2 1 4 9 7 18 9 4 1 3 0 7 6 2 2 6 0 0 1 1 3 -1 3 5 -1 8 34 1 4 2 2 7 12 3 2 0"
$ws_YDS.Range("C3").Value = $ws_YDS.Range("C3").Value2 + " 16 9 36 12 10 12 5 25 13 16 12 22 -2 2 28 9 7 7 2 1"
$ws_ST.Range("B4").Value = $ws_ST.Range("B4").Value2 + " 55 64 62"
$ws_ST.Range("B5").Value = $ws_ST.Range("B5").Value2 + " 21 5 13"
$ws_ST.Range("B6").Value = $ws_ST.Range("B6").Value2 + " 79 28 32"
$ws_ST.Range("D3").Value = $ws_ST.Range("D3").Value2 + " 33"
$ws_ST.Range("D4").Value = $ws_ST.Range("D4").Value2 + " 0"
$ws_ST.Range("D5").Value = $ws_ST.Range("D5").Value2 + " 0"

# --- Update Week 13 running totals across summary tables ---
$ws_OFF = $wb.Worksheets.Item("OFF")
$ws_DEF = $wb.Worksheets.Item("DEF")
$ws_ST = $wb.Worksheets.Item("ST")
$ws_TURNS = $wb.Worksheets.Item("TURNS")
$ws_PEN = $wb.Worksheets.Item("PEN")
$ws_OFF.Range("C2").Value = 314
$ws_OFF.Range("D2").Value = 15
$ws_OFF.Range("F2").Value = 90
$ws_OFF.Range("G2").Value = 110
$ws_OFF.Range("L2").Value = 478
$ws_OFF.Range("M2").Value = 289
$ws_OFF.Range("O2").Value = 36
$ws_OFF.Range("P2").Value = 16
$ws_OFF.Range("Q2").Value = 858
$ws_OFF.Range("C3").Value = 312
$ws_OFF.Range("E3").Value = 51
$ws_OFF.Range("F3").Value = 178
$ws_OFF.Range("I3").Value = 117
$ws_OFF.Range("J3").Value = 105
$ws_OFF.Range("N3").Value = 38
$ws_DEF.Range("C2").Value = 387
$ws_DEF.Range("D2").Value = 25
$ws_DEF.Range("E2").Value = 15
$ws_DEF.Range("F2").Value = 115
$ws_DEF.Range("G2").Value = 94
$ws_DEF.Range("I2").Value = 8
$ws_DEF.Range("J2").Value = 47
$ws_DEF.Range("L2").Value = 510
$ws_DEF.Range("M2").Value = 349
$ws_DEF.Range("O2").Value = 24
$ws_DEF.Range("P2").Value = 13
$ws_DEF.Range("Q2").Value = 928
$ws_DEF.Range("C3").Value = 300
$ws_DEF.Range("E3").Value = 65
$ws_DEF.Range("F3").Value = 187
$ws_DEF.Range("G3").Value = 65
$ws_DEF.Range("H3").Value = 45
$ws_DEF.Range("I3").Value = 103
$ws_DEF.Range("J3").Value = 96
$ws_DEF.Range("N3").Value = 33
$ws_ST.Range("B2").Value = 108
$ws_ST.Range("D2").Value = 123
$ws_ST.Range("F2").Value = 16
$ws_ST.Range("H2").Value = 11
$ws_ST.Range("B3").Value = 61
$ws_TURNS.Range("B2").Value = 18
$ws_TURNS.Range("E2").Value = 20
$ws_TURNS.Range("D3").Value = 6
$ws_TURNS.Range("E3").Value = 18
$ws_PEN.Range("D2").Value = 19
$ws_PEN.Range("D4").Value = 16
